$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Description"

# --- Row 2: F_01 ---
$ws.Range("A2").Value = "F_01"
$ws.Range("B2").Value = 'In the fairytale "Alice in DaSCHland," a curious young girl named Alice discovers a magical realm powered by DaSCH, a company that operates as a repository for research data in the humanities. This world is filled with technological wonders like talking robots, digital trees, and houses that change color. Guided by a robot named RabbIT, Alice explores various marvels, including a café with size-altering tea and a library with books that speak. She meets the Queen of Hearts, an engineer who creates devices connecting emotions with digital data. Although Alice returns to her world with a token from DaSCHland, the experience leaves her inspired, reminding her that the realm of innovation and dreams, where research data in the humanities are brought to life, awaits her return.'

# --- Row 3: F_02 ---
$ws.Range("A3").Value = "F_02"
$ws.Range("B3").Value = "
In the enchanting tale of `"Alice in DaSCHland,`" a young, adventurous girl named Alice finds herself in a fantastical world powered by DaSCH, a pioneering company known for its repository of humanities research data. DaSCHland is a place where technology and imagination merge, featuring sentient robots, trees with digital leaves, and dynamically transforming houses. With the guidance of RabbIT, a robot companion, Alice explores incredible sights, such as a café that offers teas causing her to shrink or grow and a library where books narrate their own stories. She encounters the Queen of Hearts, a brilliant engineer who designs devices that weave human emotions into digital formats. When Alice returns to her own reality, carrying a small, glowing keepsake from her journey, she is left with a sense of wonder and a knowledge that DaSCHland's visionary world, where the humanities data comes alive, will always be there, inviting her to explore further."

# --- Formatting: B2 gets font2 + vertical top + wrap ---
$ws.Range("B2").Font.Name = "Arial"
$ws.Range("B2").Font.Size = 10
$ws.Range("B2").VerticalAlignment = -4160
$ws.Range("B2").WrapText = $true

# --- Formatting: A2:A10 gets font2 + vertical top (no wrap) ---
$ws.Range("A2:A10").Font.Name = "Arial"
$ws.Range("A2:A10").Font.Size = 10
$ws.Range("A2:A10").VerticalAlignment = -4160

# --- Formatting: B3 gets font2 + wrap only ---
$ws.Range("B3").Font.Name = "Arial"
$ws.Range("B3").Font.Size = 10
$ws.Range("B3").WrapText = $true

# --- Formatting: A11:A12 gets font2 (plain) ---
$ws.Range("A11:A12").Font.Name = "Arial"
$ws.Range("A11:A12").Font.Size = 10

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 87
$ws.Rows.Item(3).RowHeight = 140
$ws.Rows.Item(4).RowHeight = 13
$ws.Rows.Item(5).RowHeight = 13
$ws.Rows.Item(6).RowHeight = 13
$ws.Rows.Item(7).RowHeight = 13
$ws.Rows.Item(8).RowHeight = 13
$ws.Rows.Item(9).RowHeight = 13
$ws.Rows.Item(10).RowHeight = 13
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 15.75

# --- Selection ---
$ws.Range("A2").Select()
